# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "61.849.71", "1.00") that must
# stay plain text, matching the workbook's original inline-string representation.
# Force text format per cell before writing so Excel does not coerce them into
# numbers, then restore the default "Normal" style so no stray number format lingers.
$priceUpdates = @{
    "D2" = '61.849.71'
    "D3" = '2.984.95'
    "D5" = '541.01'
    "D6" = '135.08'
    "D7" = '1.00'
    "D8" = '2.978.71'
    "D13" = '0.0000217'
    "D14" = '33.58'
    "D15" = '3.440.45'
    "D16" = '61.831.48'
    "D18" = '2.983.11'
    "D20" = '464.33'
    "D21" = '13.39'
    "D22" = '0.650'
    "D24" = '79.21'
    "D25" = '12.47'
    "D26" = '0.999'
    "D29" = '1.00'
    "D30" = '1.98'
    "D31" = '25.16'
    "D33" = '2.29'
    "D34" = '5.45'
    "D35" = '53.84'
    "D36" = '5.78'
    "D37" = '447.33'
    "D38" = '0.0800'
    "D39" = '0.0384'
    "D40" = '2.935.96'
    "D42" = '7.98'
    "D43" = '2.43'
    "D44" = '26.33'
    "D46" = '0.245'
    "D47" = '0.108'
    "D48" = '1.97'
    "D49" = '114.58'
}
foreach ($ref in $priceUpdates.Keys) {
    $ws.Range($ref).NumberFormat = "@"
}
foreach ($ref in $priceUpdates.Keys) {
    $ws.Range($ref).Value = $priceUpdates[$ref]
}
foreach ($ref in $priceUpdates.Keys) {
    $ws.Range($ref).Style = "Normal"
}

# Coin name / link swaps and volume(1h) percentage updates (plain text, no coercion risk)
$textUpdates = @{
    "B29" = 'FirstDigitalUSD'
    "B30" = 'ImmutableX'
    "B47" = 'Stellar'
    "B48" = 'Fetch.AI'
    "C29" = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    "C30" = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "C47" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "C48" = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    "E2" = '  -0.41%  '
    "E3" = '  -1.42%  '
    "E4" = '  +0.12%  '
    "E5" = '  -1.48%  '
    "E6" = '  -0.71%  '
    "E7" = '  +0.02%  '
    "E8" = '  -1.52%  '
    "E9" = '  -2.66%  '
    "E10" = '  +9.16%  '
    "E11" = '  -2.86%  '
    "E12" = '  -1.94%  '
    "E13" = '  -3.15%  '
    "E14" = '  -2.67%  '
    "E15" = '  -2.10%  '
    "E16" = '  -0.44%  '
    "E17" = '  -2.08%  '
    "E18" = '  -1.25%  '
    "E19" = '  -2.45%  '
    "E20" = '  -3.16%  '
    "E21" = '  +0.47%  '
    "E22" = '  -4.16%  '
    "E23" = '  +0.34%  '
    "E24" = '  -1.94%  '
    "E25" = '  +2.33%  '
    "E26" = '  -0.06%  '
    "E27" = '  -1.26%  '
    "E28" = '  -3.97%  '
    "E29" = '  +0.21%  '
    "E30" = '  +2.89%  '
    "E31" = '  -2.78%  '
    "E32" = '  -3.40%  '
    "E33" = '  -1.52%  '
    "E34" = '  -0.66%  '
    "E35" = '  -2.83%  '
    "E36" = '  -2.93%  '
    "E37" = '  -2.61%  '
    "E38" = '  -0.07%  '
    "E39" = '  -0.79%  '
    "E40" = '  -9.23%  '
    "E41" = '  -4.60%  '
    "E42" = '  -2.48%  '
    "E43" = '  -2.16%  '
    "E44" = '  +1.87%  '
    "E46" = '  -0.67%  '
    "E47" = '  -0.54%  '
    "E48" = '  -1.76%  '
    "E49" = '  -3.54%  '
    "E50" = '  -2.69%  '
    "E51" = '  -2.18%  '
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

